$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update the header timestamp / description cell ---------------------
$ws.Range("A1").Value = "Averages Set output from blixt_rp python library on 2023-06-17T15:03:30.168244"

# --- 2. Insert 5 fresh rows before row 37 ------------------------------------
# This shifts the existing "Sand D_*_sand" .. "Shale C_WELL_F_shale" block
# (old rows 37-56) down to rows 42-61, matching the new A1:BD61 dimension.
$ws.Rows("37:41").Insert()

# --- 3. Populate the 5 new rows with the "gas_sand" variant of the Sand D group --
function Set-RowValues($ws, $rowNum, $values) {
    $arr = New-Object "object[,]" 1,$values.Count
    for ($i = 0; $i -lt $values.Count; $i++) { $arr[0,$i] = $values[$i] }
    $rangeAddr = "A${rowNum}:BD${rowNum}"
    $ws.Range($rangeAddr).Value = $arr
}

$row37Values = @(
    "Sand D_gas_sand", "NONE", "MD", -999.25, -999.25, -999.25, 3336.283993999999, 1981.303252,
    2.29057, 3314.8533, 1938.96015, 2.30575, 3336.283993999999, 1981.303252, 2.29057, "NONE",
    0.141292, 0.03245786708950544, -999.25, -999.25, -999.25, -999.25, -999.25, -999.25,
    -999.25, -999.25, -999.25, -999.25, -999.25, 216.8481764696428, 125.0611344791422, 0.06275597262412556,
    -999.25, -999.25, -999.25, 0.7686486370092116, 0.2196510636464898, 0.1163338668271429, -999.25, -999.25,
    -999.25, -999.25, -999.25, -999.25, -999.25, -999.25, -999.25, -999.25,
    -999.25, 0.1642, 0.1296120796839554, "None", 0, "Volume<0.50, Porosity>0.10", "P velocity: vp_sg08, S velocity: vs_sg08, Density: rho_sg08, Porosity: phie, Volume: vcl", 45094.62744013961
)
Set-RowValues $ws 37 $row37Values

$row38Values = @(
    "Sand D_WELL_A_gas_sand", "NONE", "MD", -999.25, -999.25, -999.25, 3336.283994, 1981.303252,
    2.29057, 3314.8533, 1938.96015, 2.30575, 3336.283994, 1981.303252, 2.29057, "NONE",
    0.141292, 0.03245786708950543, -999.25, -999.25, -999.25, -999.25, -999.25, -999.25,
    -999.25, -999.25, -999.25, -999.25, -999.25, 216.8481764696428, 125.0611344791422, 0.06275597262412556,
    -999.25, -999.25, -999.25, 0.7686486370092114, 0.21965106364649, 0.116333866827143, -999.25, -999.25,
    -999.25, -999.25, -999.25, -999.25, -999.25, -999.25, -999.25, -999.25,
    -999.25, 0.1642, 0.1296120796839554, "None", 0, "Volume<0.50, Porosity>0.10", "P velocity: vp_sg08, S velocity: vs_sg08, Density: rho_sg08, Porosity: phie, Volume: vcl", 45094.62744049744
)
Set-RowValues $ws 38 $row38Values

$row39Values = @(
    "Sand D_WELL_B_gas_sand", "NONE", "MD", -999.25, -999.25, -999.25, $null, $null,
    $null, $null, $null, $null, $null, $null, $null, "NONE",
    $null, $null, -999.25, -999.25, -999.25, -999.25, -999.25, -999.25,
    -999.25, -999.25, -999.25, -999.25, -999.25, $null, $null, $null,
    -999.25, -999.25, -999.25, $null, $null, $null, -999.25, -999.25,
    -999.25, -999.25, -999.25, -999.25, -999.25, -999.25, -999.25, -999.25,
    -999.25, $null, $null, "None", 0, "Volume<0.50, Porosity>0.10", "P velocity: vp_sg08, S velocity: vs_sg08, Density: rho_sg08, Porosity: phie, Volume: vcl", 45094.62744086684
)
Set-RowValues $ws 39 $row39Values

$row40Values = @(
    "Sand D_WELL_C_gas_sand", "NONE", "MD", -999.25, -999.25, -999.25, $null, $null,
    $null, $null, $null, $null, $null, $null, $null, "NONE",
    $null, $null, -999.25, -999.25, -999.25, -999.25, -999.25, -999.25,
    -999.25, -999.25, -999.25, -999.25, -999.25, $null, $null, $null,
    -999.25, -999.25, -999.25, $null, $null, $null, -999.25, -999.25,
    -999.25, -999.25, -999.25, -999.25, -999.25, -999.25, -999.25, -999.25,
    -999.25, $null, $null, "None", 0, "Volume<0.50, Porosity>0.10", "P velocity: vp_sg08, S velocity: vs_sg08, Density: rho_sg08, Porosity: phie, Volume: vcl", 45094.62744123589
)
Set-RowValues $ws 40 $row40Values

$row41Values = @(
    "Sand D_WELL_F_gas_sand", "NONE", "MD", -999.25, -999.25, -999.25, 3336.283994, 1981.303252,
    2.29057, 3314.8533, 1938.96015, 2.30575, 3336.283994, 1981.303252, 2.29057, "NONE",
    0.141292, 0.03245786708950543, -999.25, -999.25, -999.25, -999.25, -999.25, -999.25,
    -999.25, -999.25, -999.25, -999.25, -999.25, 216.8481764696428, 125.0611344791422, 0.06275597262412556,
    -999.25, -999.25, -999.25, 0.7686486370092114, 0.21965106364649, 0.116333866827143, -999.25, -999.25,
    -999.25, -999.25, -999.25, -999.25, -999.25, -999.25, -999.25, -999.25,
    -999.25, 0.1642, 0.1296120796839554, "None", 0, "Volume<0.50, Porosity>0.10", "P velocity: vp_sg08, S velocity: vs_sg08, Density: rho_sg08, Porosity: phie, Volume: vcl", 45094.62744161212
)
Set-RowValues $ws 41 $row41Values

# --- 4. Refresh the "DateAdded" (BD) timestamp for every data row (22-61) -------
# to the new run time, mirroring the re-run of the averaging script.
$bdTimestamps = [ordered]@{
    "22" = 45094.62743249094
    "23" = 45094.6274326644
    "24" = 45094.62743286064
    "25" = 45094.62743307996
    "26" = 45094.62743329897
    "27" = 45094.62743473646
    "28" = 45094.62743497888
    "29" = 45094.6274352559
    "30" = 45094.62743552139
    "31" = 45094.62743579843
    "32" = 45094.62743731824
    "33" = 45094.6274376299
    "34" = 45094.62743794156
    "35" = 45094.62743827674
    "36" = 45094.62743862309
    "37" = 45094.62744013961
    "38" = 45094.62744049744
    "39" = 45094.62744086684
    "40" = 45094.62744123589
    "41" = 45094.62744161212
    "42" = 45094.62744320508
    "43" = 45094.6274436091
    "44" = 45094.62744519018
    "45" = 45094.62744562882
    "46" = 45094.62744609053
    "47" = 45094.62744785664
    "48" = 45094.62744834147
    "49" = 45094.62744882129
    "50" = 45094.62744929425
    "51" = 45094.62744979094
    "52" = 45094.62745156111
    "53" = 45094.62745208057
    "54" = 45094.62745258845
    "55" = 45094.62745313098
    "56" = 45094.62745367493
    "57" = 45094.62745553641
    "58" = 45094.62745610388
    "59" = 45094.62745669258
    "60" = 45094.62745730436
    "61" = 45094.62745789922
}
foreach ($r in $bdTimestamps.Keys) {
    $ws.Range("BD$r").Value = $bdTimestamps[$r]
}
